$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# ---------------------------------------------------------------------------
# New ballots: Steve Simmons (26), Carl Steward (27), Ross Newhan (28),
# Luis Rangel (29)
# ---------------------------------------------------------------------------

# Row 26 - Steve Simmons / Toronto Sun
$ws.Range("A26").Value = "Steve Simmons"
foreach ($c in @("E","I","J","K","O","Q","T","V")) {
    $ws.Range($c + "26").Value = "x"
}
$ws.Range("AK26").Value = 8
$ws.Range("AL26").Value = "Toronto Sun"
$ws.Range("AM26").Value = 43442

# Row 27 - Carl Steward / Twitter
$ws.Range("A27").Value = "Carl Steward"
foreach ($c in @("C","D","E","H","I","K","O","Q","T","V")) {
    $ws.Range($c + "27").Value = "x"
}
$ws.Range("AK27").Value = 10
$ws.Range("AL27").Value = "Twitter"
$ws.Range("AM27").Value = 43442

# Row 28 - Ross Newhan / Twitter
$ws.Range("A28").Value = "Ross Newhan"
foreach ($c in @("E","H","I","J","K","O","Q","T","V")) {
    $ws.Range($c + "28").Value = "x"
}
$ws.Range("AK28").Value = 9
$ws.Range("AL28").Value = "Twitter"
$ws.Range("AM28").Value = 43442

# Row 29 - Luis Rangel / Twitter
$ws.Range("A29").Value = "Luis Rangel"
foreach ($c in @("C","D","E","I","K","N","O","Q","R","S")) {
    $ws.Range($c + "29").Value = "x"
}
$ws.Range("AK29").Value = 10
$ws.Range("AL29").Value = "Twitter"
$ws.Range("AM29").Value = 43443

# Apply the same date style already used by the existing AM column (style
# index with numFmtId 14) to the four new date cells, reusing the format
# rather than creating a new one.
$ws.Range("AM25").Copy()
$ws.Range("AM26:AM29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the selection/viewport to match the author's final cursor position.
[void]$ws.Range("G28").Select()

Write-Output "done"
